$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 13151
$ws1.Range("F5").Value  = 1
$ws1.Range("F6").Value  = 100
$ws1.Range("F7").Value  = 108
$ws1.Range("F9").Value  = 35
$ws1.Range("F10").Value = 22
$ws1.Range("F11").Value = 13109
$ws1.Range("F12").Value = 315
$ws1.Range("F13").Value = 559
$ws1.Range("F14").Value = 8802
$ws1.Range("F15").Value = 7864
$ws1.Range("F16").Value = 221
$ws1.Range("F27").Value = 55

# Sheet "全部类型" (All Types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 13151
$ws4.Range("F6").Value  = 1
$ws4.Range("F7").Value  = 100
$ws4.Range("F8").Value  = 108
$ws4.Range("F10").Value = 35
$ws4.Range("F11").Value = 22
$ws4.Range("F12").Value = 13109
$ws4.Range("F13").Value = 315
$ws4.Range("F14").Value = 559
$ws4.Range("F15").Value = 8802
$ws4.Range("F16").Value = 7864
$ws4.Range("F17").Value = 221
$ws4.Range("F30").Value = 55
